# Insert a new weekly data row for "Hortaliza, Terminal La Palmera de La Serena - Jengibre"
# before the existing row 26, shifting the following rows (old 26-54) down to (27-55).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26 (pushes rows 26..54 down to 27..55)
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the new weekly record.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R are constant across this subset,
# mirror them from the row above (row 25).
$ws.Cells.Item(26, 1).Value = $ws.Cells.Item(25, 1).Value()
$ws.Cells.Item(26, 2).Value = $ws.Cells.Item(25, 2).Value()
$ws.Cells.Item(26, 3).Value = $ws.Cells.Item(25, 3).Value()

$ws.Cells.Item(26, 4).Value = 44763

$ws.Cells.Item(26, 5).Value = $ws.Cells.Item(25, 5).Value()
$ws.Cells.Item(26, 6).Value = $ws.Cells.Item(25, 6).Value()
$ws.Cells.Item(26, 7).Value = $ws.Cells.Item(25, 7).Value()
$ws.Cells.Item(26, 8).Value = $ws.Cells.Item(25, 8).Value()
$ws.Cells.Item(26, 9).Value = $ws.Cells.Item(25, 9).Value()

$ws.Cells.Item(26, 10).Value = 500
$ws.Cells.Item(26, 11).Value = 15000
$ws.Cells.Item(26, 12).Value = 16000
$ws.Cells.Item(26, 13).Value = 15500

$ws.Cells.Item(26, 14).Value = $ws.Cells.Item(25, 14).Value()
$ws.Cells.Item(26, 15).Value = $ws.Cells.Item(25, 15).Value()

$ws.Cells.Item(26, 16).Value = 1192

$ws.Cells.Item(26, 17).Value = $ws.Cells.Item(25, 17).Value()
$ws.Cells.Item(26, 18).Value = $ws.Cells.Item(25, 18).Value()

# Make sure the D column keeps the date number format used by the rest of the column.
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(25, 4).NumberFormat()
